# Split the single run of concatenated bibliography references into
# separate entries divided by a blank line (two manual line breaks),
# matching the authored diff that inserted <w:br/><w:br/> between each
# reference.

$d = $word.ActiveDocument

$pairs = @(
    @{ Find = "ra Atlas, 2005.JAFFE, R. W. Ad"; Replace = "ra Atlas, 2005.^l^lJAFFE, R. W. Ad" },
    @{ Find = "ra Atlas, 2002.GITMAN, L. J. -"; Replace = "ra Atlas, 2002.^l^lGITMAN, L. J. -" },
    @{ Find = " Pearson, 2017.ASSAF NETO, A. "; Replace = " Pearson, 2017.^l^lASSAF NETO, A. " },
    @{ Find = "lo: Atlas, 2014MORANTE, A. S. "; Replace = "lo: Atlas, 2014^l^lMORANTE, A. S. " },
    @{ Find = "o: Atlas, 2009.NEWNAN, D. G.; "; Replace = "o: Atlas, 2009.^l^lNEWNAN, D. G.; " },
    @{ Find = "ulo: LTC, 2000.HOJI, M.; LUZ, "; Replace = "ulo: LTC, 2000.^l^lHOJI, M.; LUZ, " },
    @{ Find = "o: Atlas, 2019.SANVICENTE, A. "; Replace = "o: Atlas, 2019.^l^lSANVICENTE, A. " }
)

foreach ($pair in $pairs) {
    $range = $d.Content
    $ok = $range.Find.Execute($pair.Find, $true, $false, $false, $false, $false, $true, 1, $false, $pair.Replace, 2)
    Write-Output "Replaced '$($pair.Find)' -> $ok"
}
